$d = $word.ActiveDocument

# Replaces the text of every paragraph whose (paragraph-mark-trimmed) text
# exactly equals $findText with $newText.
#
# We use InsertXML (instead of Range.Text = ...) because this engine
# coalesces / rebuilds a paragraph's runs whenever Range.Text is assigned,
# which would silently drop sibling empty runs (<w:r/>) that are present
# in the original document and must be kept unchanged. InsertXML on a
# Range constructed via $d.Range(start, end) replaces only that range's
# content and leaves the rest of the paragraph (pPr, sibling runs) intact.
#
# We detect the exact run-level formatting (w:rPr, e.g. bold/italic) of
# the run actually being replaced by reading the paragraph's own raw
# OOXML (Paragraph.Range.WordOpenXML, which - unlike a range built with
# $d.Range(...) - correctly reports only that paragraph's XML) instead of
# relying on the computed/effective Range.Bold / Range.Italic properties,
# since those reflect inherited style formatting too (e.g. Heading1 is
# bold via its style, even though the run itself carries no w:rPr).
function Replace-ParagraphText {
    param(
        [string]$findText,
        [string]$newText
    )

    $wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs.Item($i)
        $paraText = $para.Range.Text.TrimEnd([char]13)

        if ($paraText -eq $findText) {
            $rawXml = $para.Range.WordOpenXML
            $bodyIdx = $rawXml.IndexOf("<w:body")
            $endIdx = $rawXml.IndexOf("</w:p>", $bodyIdx)
            $paraXml = $rawXml.Substring($bodyIdx, $endIdx - $bodyIdx)

            $rPrXml = ""
            if ($paraXml -match '<w:rPr>.*?</w:rPr>') {
                $rPrXml = $matches[0]
            }

            $escaped = $newText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

            $rng = $d.Range($para.Range.Start, $para.Range.End)
            $xml = '<w:p ' + $wNs + '><w:r>' + $rPrXml + '<w:t>' + $escaped + '</w:t></w:r></w:p>'
            $rng.InsertXML($xml)
        }
    }
}

Replace-ParagraphText "Play Manhattan Goes Wild Free - Read Our Game Review Now!" "Play Manhattan Goes Wild for Free - Review & Gameplay"
Replace-ParagraphText "Infectious wild symbol that morphs other symbols into wilds" "High-volatility gameplay with 243 ways to win"
Replace-ParagraphText "Bonus feature offers free spins and an extra-spin meter" "Customizable Autoplay feature"
Replace-ParagraphText "Wide range of bets from 0.10 to 100 per spin" "Wide range of betting options suitable for all players"
Replace-ParagraphText "High RTP of 96.16%" "Unique wild symbol and bonus features"
Replace-ParagraphText "Paylines are not customizable" "Limited options for customizing paylines"
Replace-ParagraphText "Wondering if Manhattan Goes Wild is worth playing for free? Read our review now and discover the game's unique features and winning potential." "Discover the high-volatility slot game Manhattan Goes Wild with unique wild symbols and bonus features. Play for free now!"
